$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12/13 and 49/50 swap Coin/Link (B/C) together with Price/Volume (D/E).
# All other rows only update Price (D) and Volume(1h) (E).

# Row 2
$ws.Cells.Item(2, 4).Value = "29.022.81"
$ws.Cells.Item(2, 5).Value = "  -2.20%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.987.49"
$ws.Cells.Item(3, 5).Value = "  -1.44%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.015"
$ws.Cells.Item(4, 5).Value = "  +0.26%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "329.79"
$ws.Cells.Item(5, 5).Value = "  -0.91%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.21%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4941"
$ws.Cells.Item(7, 5).Value = "  -2.40%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4169"
$ws.Cells.Item(8, 5).Value = "  -2.24%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "55.26"
$ws.Cells.Item(9, 5).Value = "  +2.06%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08826"
$ws.Cells.Item(10, 5).Value = "  -4.58%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.084"
$ws.Cells.Item(11, 5).Value = "  -4.05%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "2.053.73"
$ws.Cells.Item(12, 5).Value = "  +5.60%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "Solana"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "22.80"
$ws.Cells.Item(13, 5).Value = "  -3.55%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.928"
$ws.Cells.Item(14, 5).Value = "  -2.54%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.378"
$ws.Cells.Item(15, 5).Value = "  -2.99%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.015"
$ws.Cells.Item(16, 5).Value = "  +0.31%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "91.71"
$ws.Cells.Item(17, 5).Value = "  -4.18%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001099"
$ws.Cells.Item(18, 5).Value = "  -2.55%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06665"
$ws.Cells.Item(19, 5).Value = "  -0.14%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "19.34"
$ws.Cells.Item(20, 5).Value = "  -3.33%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.014"
$ws.Cells.Item(21, 5).Value = "  +0.54%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.951"
$ws.Cells.Item(22, 5).Value = "  -1.00%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "29.080.79"
$ws.Cells.Item(23, 5).Value = "  -2.12%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.90"
$ws.Cells.Item(24, 5).Value = "  -1.30%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.312"
$ws.Cells.Item(25, 5).Value = "  +1.81%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.275.79"
$ws.Cells.Item(26, 5).Value = "  +3.74%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.73"
$ws.Cells.Item(27, 5).Value = "  -0.62%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "156.71"
$ws.Cells.Item(28, 5).Value = "  -1.98%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "6.207"
$ws.Cells.Item(29, 5).Value = "  -3.35%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.225"
$ws.Cells.Item(30, 5).Value = "  -5.61%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "126.23"
$ws.Cells.Item(31, 5).Value = "  -2.03%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.037"
$ws.Cells.Item(32, 5).Value = "  -2.51%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.09849"
$ws.Cells.Item(33, 5).Value = "  -1.54%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.525"
$ws.Cells.Item(34, 5).Value = "  -4.77%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.817"
$ws.Cells.Item(35, 5).Value = "  -1.31%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.734"
$ws.Cells.Item(36, 5).Value = "  -1.85%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.02402"
$ws.Cells.Item(37, 5).Value = "  -3.01%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.306"
$ws.Cells.Item(38, 5).Value = "  -1.65%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "9.015"
$ws.Cells.Item(39, 5).Value = "  -6.37%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06358"
$ws.Cells.Item(40, 5).Value = "  -0.93%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.6440"
$ws.Cells.Item(41, 5).Value = "  -2.50%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "11.46"
$ws.Cells.Item(42, 5).Value = "  -3.35%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1970"
$ws.Cells.Item(43, 5).Value = "  -5.85%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.012"
$ws.Cells.Item(44, 5).Value = "  +0.25%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +4.76%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.6146"
$ws.Cells.Item(46, 5).Value = "  -3.86%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "13.25"
$ws.Cells.Item(47, 5).Value = "  -2.81%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.151"
$ws.Cells.Item(48, 5).Value = "  -3.31%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "PancakeSwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.486"
$ws.Cells.Item(49, 5).Value = "  -1.53%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00000000343"
$ws.Cells.Item(50, 5).Value = "  +5.46%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.152"
$ws.Cells.Item(51, 5).Value = "  +6.02%  "
